# INDEX can now work with ranges passed as its second and/or third arguments.
# Demonstrate this on Sheet1 by rebuilding the SUMIF helper column as a proper
# shared formula and adding a new example that uses ranges for INDEX's
# row/column arguments.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Re-enter the SUMIF formula across G43:G48 in one shot so it is stored as a
# single shared formula (matching the refactor in the sample workbook).
$ws.Range("G43:G48").Formula = "=SUMIF(`$E`$43:`$E`$48,`$F`$43:`$F`$48,`$D`$43:`$D`$48)"

# New example data for the INDEX-with-ranges demonstration.
$ws.Range("D53").Value = 10
$ws.Range("E53").Value = 20
$ws.Range("F53").Value = 30
$ws.Range("G53").Value = 40
$ws.Range("H53").Value = 50

$ws.Range("D54").Value = 5
$ws.Range("E54").Value = 4
$ws.Range("F54").Value = 3
$ws.Range("G54").Value = 2
$ws.Range("H54").Value = 1

# INDEX called with ranges (rather than single cells) for its row/column
# arguments - entered across the row as one shared formula.
$ws.Range("D55:H55").Formula = "=INDEX(`$D`$53:`$H`$53,1,`$D`$54:`$H`$54)"

$ws.Range("H55").Select()
